# Commit: "changing paths for new repo structure"
#
# 1. Repoint the catalysts-data CSV reference (used throughout the
#    "Source" sheet) from the old mappings-relative path to the new
#    processed-data path under the repo's new layout.
# 2. Leave the "Source" sheet scrolled/selected further up (closer to
#    where the edit was made) instead of all the way at the bottom.

$wb = $excel.ActiveWorkbook

$source = $wb.Worksheets.Item("Source")
$source.Activate() | Out-Null

# Update every cell that held the old relative path (column C, all data
# rows) to the new path. Replace() rewrites the shared string in place so
# every cell referencing it is updated together.
$usedRange = $source.Range("A1:C113")
$usedRange.Replace("../mappings/data/catalystsdata.csv", "../data/processed/catalystsdata.csv") | Out-Null

# Move the viewport/selection on the Source sheet.
$source.Range("A70").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 70
$excel.ActiveWindow.ScrollColumn = 1
$source.Range("C108").Select() | Out-Null
